# DDAf_2023_Tableau_annexe_Tab12.xlsx -- "Add files via upload" update
#
# Three country-label edits (adds/removes the "resource-rich" footnote
# marker "*") plus the knock-on recompute of the regional/group summary
# rows (68-98) whose figures are pasted-in aggregates (no live formulas
# in this sheet) that shift once South Sudan flips into, and Nigeria
# flips out of, the "resource-rich" bucket.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab12")

# ---------------------------------------------------------------------
# 1) Country-name edits
# ---------------------------------------------------------------------
# South Sudan becomes a resource-rich country -> add the "*" marker
$ws.Range("B34").Value = "Soudan du Sud*"
# Cabo Verde becomes a resource-rich country -> add the "*" marker
$ws.Range("B48").Value = "Cabo Verde*"
# Nigeria is no longer flagged resource-rich -> drop the "*" marker
$ws.Range("B57").Value = "Nigeria"

# ---------------------------------------------------------------------
# 2) Matching shading: resource-rich country rows are highlighted with a
#    light-blue fill. Re-use the formatting of an existing correctly
#    shaded / unshaded row so fonts, borders and number formats stay
#    consistent with the rest of the table.
# ---------------------------------------------------------------------
# Soudan du Sud (row 34) gains the "resource-rich" shading (copy format
# from Tchad*, row 17, which already carries that exact style).
$ws.Range("B17:O17").Copy() | Out-Null
$ws.Range("B34:O34").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Nigeria (row 57) loses the "resource-rich" shading (copy format from
# Somalie, row 33, an ordinary unshaded row).
$ws.Range("B33:O33").Copy() | Out-Null
$ws.Range("B57:O57").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Recomputed aggregate rows (static pasted-in values, not formulas)
# ---------------------------------------------------------------------
$ws.Range("E68").Value = 25.4444444444444
$ws.Range("I68").Value = 53.941176470588303
$ws.Range("C69").Value = 64.616666666666703
$ws.Range("D69").Value = 45.35
$ws.Range("E69").Value = 19.5833333333333
$ws.Range("F69").Value = 26.216666666666701
$ws.Range("G69").Value = 84.4166666666667
$ws.Range("H69").Value = 69.866666666666703
$ws.Range("I69").Value = 40.549999999999997
$ws.Range("J69").Value = 48.633333333333297
$ws.Range("K69").Value = 701.91200000000003
$ws.Range("L69").Value = 2.7845174928333298
$ws.Range("M69").Value = 2135.1089999999999
$ws.Range("N69").Value = 24.844539699666701
$ws.Range("O69").Value = 28.316666666666698
$ws.Range("G72").Value = 83.183333333333294
$ws.Range("O77").Value = 52.977777777777803
$ws.Range("C80").Value = 69.066666666666706
$ws.Range("D80").Value = 42.3333333333333
$ws.Range("E80").Value = 16.8
$ws.Range("F80").Value = 25.233333333333299
$ws.Range("G80").Value = 86.933333333333394
$ws.Range("H80").Value = 75.7
$ws.Range("I80").Value = 45.5
$ws.Range("J80").Value = 52.8
$ws.Range("K80").Value = 132.12
$ws.Range("L80").Value = 2.9478269036000002
$ws.Range("M80").Value = 2430.0010000000002
$ws.Range("N80").Value = 39.175325964999999
$ws.Range("O80").Value = 28.7
$ws.Range("C82").Value = 69.103125000000006
$ws.Range("D82").Value = 47.34375
$ws.Range("E82").Value = 25.034375000000001
$ws.Range("F82").Value = 33.924999999999997
$ws.Range("G82").Value = 89.3
$ws.Range("H82").Value = 76.453333333333404
$ws.Range("I82").Value = 52.85
$ws.Range("J82").Value = 61.3466666666667
$ws.Range("K82").Value = 5233.5878106
$ws.Range("L82").Value = 8.54717057525
$ws.Range("M82").Value = 25186.054678699998
$ws.Range("N82").Value = 25.095729166135101
$ws.Range("O82").Value = 29.967567567567599
$ws.Range("C84").Value = 63.461111111111101
$ws.Range("D84").Value = 40.677777777777798
$ws.Range("E84").Value = 20.011111111111099
$ws.Range("F84").Value = 27.672222222222199
$ws.Range("G84").Value = 86.293750000000003
$ws.Range("H84").Value = 72.53125
$ws.Range("I84").Value = 47.762500000000003
$ws.Range("J84").Value = 55.84375
$ws.Range("K84").Value = 815.27690299999995
$ws.Range("L84").Value = 11.593017487588201
$ws.Range("M84").Value = 2119.1882406999998
$ws.Range("N84").Value = 31.246143113611101
$ws.Range("O84").Value = 19
$ws.Range("C86").Value = 73.55
$ws.Range("D86").Value = 53.131250000000001
$ws.Range("E86").Value = 28.712499999999999
$ws.Range("F86").Value = 38.987499999999997
$ws.Range("G86").Value = 91.193749999999994
$ws.Range("H86").Value = 79.1875
$ws.Range("I86").Value = 55.256250000000001
$ws.Range("J86").Value = 64.2
$ws.Range("K86").Value = 4364.6499076
$ws.Range("L86").Value = 5.75767284355
$ws.Range("M86").Value = 20321.012438000002
$ws.Range("N86").Value = 23.537399084315801
$ws.Range("O86").Value = 34.125
$ws.Range("C87").Value = 78.661538461538498
$ws.Range("D87").Value = 53.515384615384598
$ws.Range("E87").Value = 30.207692307692302
$ws.Range("F87").Value = 40.815384615384602
$ws.Range("G87").Value = 93.915384615384596
$ws.Range("H87").Value = 77.615384615384599
$ws.Range("I87").Value = 57.069230769230799
$ws.Range("J87").Value = 66.5230769230769
$ws.Range("K87").Value = 81611.359931700004
$ws.Range("L87").Value = 7.5778063753846201
$ws.Range("M87").Value = 203608.65209049999
$ws.Range("N87").Value = 21.662767979689701
$ws.Range("O87").Value = 43.466666666666697
$ws.Range("C89").Value = 83.507692307692295
$ws.Range("D89").Value = 70.730769230769198
$ws.Range("E89").Value = 52.307692307692299
$ws.Range("F89").Value = 60.546153846153899
$ws.Range("G89").Value = 93.224999999999994
$ws.Range("H89").Value = 91.575000000000003
$ws.Range("I89").Value = 78.45
$ws.Range("J89").Value = 83.858333333333405
$ws.Range("K89").Value = 78385.100199499997
$ws.Range("L89").Value = 6.64171165322223
$ws.Range("M89").Value = 258205.04571020001
$ws.Range("N89").Value = 22.0593773816341
$ws.Range("O89").Value = 60.129032258064498
$ws.Range("K90").Value = 450359.02202899999
$ws.Range("L90").Value = 11.3053117095814
$ws.Range("M90").Value = 2654923.2574481
$ws.Range("N90").Value = 39.405172897615401
$ws.Range("O90").Value = 82.027083333333394
$ws.Range("O94").Value = 55.4
$ws.Range("L95").Value = 8.9169743408333293
$ws.Range("C97").Value = 65.707142857142898
$ws.Range("D97").Value = 42.717857142857198
$ws.Range("E97").Value = 21.05
$ws.Range("F97").Value = 29.35
$ws.Range("G97").Value = 87.126923076923106
$ws.Range("H97").Value = 72.723076923076903
$ws.Range("I97").Value = 47.938461538461503
$ws.Range("J97").Value = 56.342307692307699
$ws.Range("K97").Value = 2096.4389030000002
$ws.Range("L97").Value = 9.3205861987407399
$ws.Range("M97").Value = 5913.8902406999996
$ws.Range("N97").Value = 27.161026947185199
$ws.Range("O97").Value = 23.990322580645199
$ws.Range("K98").Value = 3374.2669999999998
$ws.Range("L98").Value = 6.9485901003571504
$ws.Range("M98").Value = 10870.9432445
$ws.Range("N98").Value = 22.426116987266699
